# Auto-generated script to update cryptos list data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.931.51"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "1.828.17"
$ws.Range("E3").Value = "  -2.27%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "240.35"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6854"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.09%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07623"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.31%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3016"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -4.53%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "23.46"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.99%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07745"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "1.834.75"
$ws.Range("E12").Value = "  -2.97%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.047"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.32%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "90.27"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.09%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6733"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.55%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.456"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008267"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "28.942.81"
$ws.Range("E18").Value = "  -1.85%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "243.39"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.32%  "
$ws.Range("D20").Value = "2.099.47"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -4.22%  "
$ws.Range("E22").Value = "  +0.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.420"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("E24").Value = "  +0.00%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1473"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.58%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "161.25"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.724"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.89%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.15"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.72%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.536"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.17%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.207"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.160"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.191"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.54%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05127"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.74%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7650"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.87%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.818"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.31%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.147"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.44%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.702"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01833"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "1.223.26"
$ws.Range("E39").Value = "  -3.41%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.704"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9116"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.08%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "108.71"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "1.999.27"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5172"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "5.376"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -10.03%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.487"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.00000000120"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -8.10%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "63.10"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -12.16%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.725"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.00%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.891"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.65%  "
